$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 9.82461274400216
$ws.Range("G2").Value = 9.6226578296356
$ws.Range("H2").Value = 10.03087779307675
$ws.Range("I2").Value = 0.002528852152208033
$ws.Range("J2").Value = 0.002253088948055081
$ws.Range("K2").Value = 0.002861120764638966
$ws.Range("L2").Value = 0.008890180060885326
$ws.Range("M2").Value = 0.008697326742682056
$ws.Range("N2").Value = 0.009091865758660951
$ws.Range("F3").Value = 0.04664534464428318
$ws.Range("G3").Value = 0.04634631039237534
$ws.Range("H3").Value = 0.0469325087099207
$ws.Range("I3").Value = 0.04508553840610342
$ws.Range("J3").Value = 0.04479752780671629
$ws.Range("K3").Value = 0.04536103108637986
$ws.Range("L3").Value = 0.04669584838478093
$ws.Range("M3").Value = 0.04639689118820508
$ws.Range("N3").Value = 0.04698308668758529
$ws.Range("F4").Value = 9.871258088646442
$ws.Range("G4").Value = 9.669004140027974
$ws.Range("H4").Value = 10.07781030178667
$ws.Range("I4").Value = 0.04761439055831146
$ws.Range("J4").Value = 0.04705061675477137
$ws.Range("K4").Value = 0.04822215185101883
$ws.Range("L4").Value = 0.05558602844566626
$ws.Range("M4").Value = 0.05509421793088713
$ws.Range("N4").Value = 0.05607495244624622
